$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format first so numeric-looking strings
# (e.g. "306.76") are written back as text, matching the source
# inlineStr cells rather than being auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "42.951.76"
$ws.Range("E2").Value = "  -1.23%  "
$ws.Range("D3").Value = "2.341.79"
$ws.Range("E3").Value = "  +1.37%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "306.76"
$ws.Range("E5").Value = "  -1.32%  "
$ws.Range("D6").Value = "100.62"
$ws.Range("E6").Value = "  -1.18%  "
$ws.Range("E7").Value = "  -4.91%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -3.39%  "
$ws.Range("D10").Value = "35.07"
$ws.Range("E10").Value = "  -1.93%  "
$ws.Range("D11").Value = "52.16"
$ws.Range("E11").Value = "  +0.73%  "
$ws.Range("D12").Value = "0.0801"
$ws.Range("E12").Value = "  -1.62%  "
$ws.Range("E13").Value = "  -0.18%  "
$ws.Range("E14").Value = "  -2.70%  "
$ws.Range("D15").Value = "15.94"
$ws.Range("E15").Value = "  +6.61%  "
$ws.Range("D16").Value = "2.398.49"
$ws.Range("E16").Value = "  +3.86%  "
$ws.Range("D17").Value = "0.806"
$ws.Range("E17").Value = "  -0.22%  "
$ws.Range("D18").Value = "42.875.01"
$ws.Range("E18").Value = "  -1.17%  "
$ws.Range("D19").Value = "6.23"
$ws.Range("E19").Value = "  +0.99%  "
$ws.Range("D20").Value = "0.0₃0913"
$ws.Range("E20").Value = "  -1.66%  "
$ws.Range("D21").Value = "11.73"
$ws.Range("E21").Value = "  -4.96%  "
$ws.Range("D22").Value = "67.97"
$ws.Range("E22").Value = "  -0.21%  "
$ws.Range("D23").Value = "236.84"
$ws.Range("E23").Value = "  -1.97%  "
$ws.Range("D24").Value = "2.03"
$ws.Range("E24").Value = "  +0.43%  "
$ws.Range("E25").Value = "  -2.09%  "
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("D27").Value = "25.47"
$ws.Range("E27").Value = "  +3.42%  "
$ws.Range("E28").Value = "  +9.64%  "
$ws.Range("D29").Value = "35.13"
$ws.Range("E29").Value = "  -4.41%  "
$ws.Range("E30").Value = "  -3.15%  "
$ws.Range("D31").Value = "159.94"
$ws.Range("E31").Value = "  -4.45%  "
$ws.Range("E32").Value = "  -0.06%  "
$ws.Range("E33").Value = "  -2.71%  "
$ws.Range("D34").Value = "4.71"
$ws.Range("E34").Value = "  +8.68%  "
$ws.Range("E35").Value = "  -0.66%  "
$ws.Range("D36").Value = "17.43"
$ws.Range("E36").Value = "  -0.71%  "
$ws.Range("E37").Value = "  -2.24%  "
$ws.Range("D38").Value = "2.96"
$ws.Range("E38").Value = "  -3.75%  "
$ws.Range("E39").Value = "  -0.45%  "
$ws.Range("E40").Value = "  -3.07%  "
$ws.Range("E41").Value = "  -2.58%  "
$ws.Range("D42").Value = "2.42"
$ws.Range("E42").Value = "  +4.63%  "
$ws.Range("D43").Value = "2.021.94"
$ws.Range("E43").Value = "  +2.82%  "
$ws.Range("E44").Value = "  -1.38%  "
$ws.Range("D45").Value = "18.96"
$ws.Range("E45").Value = "  -1.37%  "
$ws.Range("D46").Value = "10.36"
$ws.Range("E46").Value = "  +3.60%  "
$ws.Range("E47").Value = "  -1.81%  "
$ws.Range("D48").Value = "56.02"
$ws.Range("E48").Value = "  +0.70%  "
$ws.Range("D49").Value = "2.90"
$ws.Range("E49").Value = "  -0.60%  "
$ws.Range("D50").Value = "2.568.29"
$ws.Range("E50").Value = "  +1.22%  "
$ws.Range("E51").Value = "  +1.65%  "
